$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 7645
    $ws.Range("F10").Value = 154
}
